# Refresh nightfall data tables:
#  - add new column T ("uint" / "xp") with xp values for existing rows 6-9
#  - append two new enemy rows (10, 11), including the new xp column
#
# Every column in this sheet stores its data as text (even numeric-looking
# values like "40" or "1.80") -- that's how the existing A4:S9 block is
# encoded, and it's flagged with ignoredError numberStoredAsText="1". A
# leading apostrophe is Excel's standard "force text" input marker, so we
# use it for every write here to land on the same Text cell type without
# touching any cell's number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column T: header + field name ---
$ws.Range("T4").Value = "'uint"
$ws.Range("T5").Value = "'xp"

# --- New xp values for the existing enemies (rows 6-9) ---
$ws.Range("T6").Value  = "'18"
$ws.Range("T7").Value  = "'26"
$ws.Range("T8").Value  = "'32"
$ws.Range("T9").Value  = "'38"

# --- Row 10: Harbor Dredger ---
$ws.Range("A10").Value = "'40"
$ws.Range("B10").Value = "'06"
$ws.Range("C10").Value = "'0005"
$ws.Range("D10").Value = "'Harbor Dredger"
$ws.Range("E10").Value = "'CONSTRUCT"
$ws.Range("F10").Value = "'520"
$ws.Range("G10").Value = "'55"
$ws.Range("H10").Value = "'2.2"
$ws.Range("I10").Value = "'MANUAL"
$ws.Range("J10").Value = "'1.10"
$ws.Range("K10").Value = "'0"
$ws.Range("L10").Value = "'0.00"
$ws.Range("M10").Value = "'fx/projectiles/dredger_slam.png"
$ws.Range("N10").Value = "'fx/impact/dredger_slam.png"
$ws.Range("O10").Value = "'FIRE"
$ws.Range("P10").Value = "'KINETIC"
$ws.Range("Q10").Value = "'loot:dredger_core"
$ws.Range("R10").Value = "'14"
$ws.Range("S10").Value = "'Heavy mech charges the beacon and causes shock tremors."
$ws.Range("T10").Value = "'46"

# --- Row 11: Myriad Fragment ---
$ws.Range("A11").Value = "'40"
$ws.Range("B11").Value = "'06"
$ws.Range("C11").Value = "'0006"
$ws.Range("D11").Value = "'Myriad Fragment"
$ws.Range("E11").Value = "'ABERRATION"
$ws.Range("F11").Value = "'160"
$ws.Range("G11").Value = "'18"
$ws.Range("H11").Value = "'4.8"
$ws.Range("I11").Value = "'BURST"
$ws.Range("J11").Value = "'1.90"
$ws.Range("K11").Value = "'26"
$ws.Range("L11").Value = "'0.75"
$ws.Range("M11").Value = "'fx/projectiles/fragment_dart.png"
$ws.Range("N11").Value = "'fx/impact/fragment_spark.png"
$ws.Range("O11").Value = "'LIGHT"
$ws.Range("P11").Value = "'VOID"
$ws.Range("Q11").Value = "'loot:fragment_cache"
$ws.Range("R11").Value = "'8"
$ws.Range("S11").Value = "'Skittering shards fire dart volleys in packs."
$ws.Range("T11").Value = "'20"
